$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear values in F1:H1 while keeping their style
$ws.Range("F1:H1").ClearContents()

# Fill C2:E2 with the same value as B2 ("Абакан")
$ws.Range("C2").Value = "Абакан"
$ws.Range("D2").Value = "Абакан"
$ws.Range("E2").Value = "Абакан"

# Clear values in A76:A78 while keeping their style
$ws.Range("A76:A78").ClearContents()

# Update the view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("D79").Select()
